# Iteration Burndown.xlsx - update burndown data for "Iteration 5" and
# refresh the chart title to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Hours Left" (column B) burndown values ---------------------
$ws.Range("B4").Value  = 11
$ws.Range("B5").Value  = 11
$ws.Range("B6").Value  = 11
$ws.Range("B7").Value  = 8
$ws.Range("B8").Value  = 8
$ws.Range("B9").Value  = 8
$ws.Range("B10").Value = 8
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 4
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 0

# --- Update the burndown chart title -------------------------------------
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.ChartTitle.Text = "Iteration 5 Burndown"

# --- Match the author's final selection -----------------------------------
$ws.Range("B15").Select()
